# WeatherData.xlsx edit: remove the first sounding-source row (the
# rucsoundings.noaa.gov link) from Sheet1, so the remaining
# weather.uwyo.edu link shifts up from A2 into A1.
#
# Resulting state: a single populated cell, A1, holding the
# weather.uwyo.edu URL text but still carrying the original A1
# hyperlink relationship (rId1, which targets rucsoundings.noaa.gov) -
# i.e. exactly what a "delete row 1 / shift up" operation produces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear every hyperlink on the sheet first. The engine does not reliably
# retarget/garbage-collect individual Hyperlinks objects once the rows
# they were anchored to get shifted by a structural delete, so the
# reliable path is: drop all hyperlinks up front, do the structural
# edit, then re-create only the single hyperlink that should remain.
$ws.Range("A1:A2").Hyperlinks.Delete()

# Delete row 1 entirely (the rucsoundings.noaa.gov row). This shifts the
# former A2 (weather.uwyo.edu) up to A1, drops the sharedString that was
# only used by the old A1, and shrinks the sheet dimension down to A1.
$ws.Rows(1).Delete()

# Re-attach the hyperlink relationship to the surviving cell so a link
# keeps working (matches the single <hyperlink ref="A1" r:id="rId1".../>
# left in the saved file).
$ws.Hyperlinks.Add($ws.Range("A1"), "https://rucsoundings.noaa.gov/get_soundings.cgi?data_source=GFS&latest=latest&start_year=2024&start_month_name=Sep&start_mday=11&start_hour=1&start_min=0&n_hrs=24&fcst_len=shortest&airport=&start=latest")

# Adding the hyperlink can nudge the cell style; make sure A1 keeps the
# workbook's existing "Hyperlink" cell style (same as before the edit).
$ws.Range("A1").Style = "Hyperlink"

# Re-select A1 (the only remaining cell) so the sheet's active selection
# no longer references the now-deleted A2.
$ws.Range("A1").Select()
